$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44630
$ws.Range("J2").Value = 90
$ws.Range("K2").Value = 2500
$ws.Range("L2").Value = 3000
$ws.Range("M2").Value = 2722
$ws.Range("O2").Value = "Región Metropolitana"
$ws.Range("P2").Value = 454

# Row 3
$ws.Range("D3").Value = 44637
$ws.Range("J3").Value = 170
$ws.Range("K3").Value = 2800
$ws.Range("L3").Value = 3000
$ws.Range("M3").Value = 2906
$ws.Range("O3").Value = "Región Metropolitana"
$ws.Range("P3").Value = 484

# Row 4
$ws.Range("D4").Value = 44643
$ws.Range("J4").Value = 90
$ws.Range("K4").Value = 2800
$ws.Range("L4").Value = 3000
$ws.Range("M4").Value = 2911
$ws.Range("O4").Value = "Región Metropolitana"
$ws.Range("P4").Value = 485

# Row 5
$ws.Range("D5").Value = 44631
$ws.Range("J5").Value = 110
$ws.Range("K5").Value = 3000
$ws.Range("L5").Value = 3500
$ws.Range("M5").Value = 3273
$ws.Range("O5").Value = "Provincia de Chacabuco"
$ws.Range("P5").Value = 546

# Row 6
$ws.Range("D6").Value = 44650
$ws.Range("J6").Value = 130
$ws.Range("K6").Value = 3000
$ws.Range("L6").Value = 3500
$ws.Range("M6").Value = 3308
$ws.Range("O6").Value = "Región Metropolitana"
$ws.Range("P6").Value = 551
